# Update column G ("K") values in the strike-count table on Sheet1.
# The data is recomputed ("regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals") and the resulting K values for
# rows 2-22 change as follows (row 18 is unchanged and left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 0
    6  = 3
    7  = 0
    8  = 0
    9  = 0
    10 = 2
    11 = 0
    12 = 2
    13 = 3
    14 = 2
    15 = 2
    16 = 1
    17 = 3
    19 = 1
    20 = 1
    21 = 0
    22 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
